$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2015, 0, 512912.8143242391, 7307.843315566571, 243.8860905380762, 661.8095107925362, 1270),
    @(2015, 0, 512912.8143242391, 7307.843315566571, 243.8860905380762, 661.8095107925362, 1270),
    @(2015, 0, 517027.4116409178, 7328.58410374858, 244.3627850925517, 667.4420417228844, 1270),
    @(2015, 0, 517027.4116409178, 7328.58410374858, 244.3627850925517, 667.4420417228844, 1270)
)

$startRow = 6
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($col = 1; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $data[$i][$col - 1]
    }
}
